$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "nome_dono" column (old column D) entirely - its data (owner name)
# is no longer tracked. This shifts observacoes/valor_compra/status/Data de Cadastro
# left by one column (D..G) and a new "Valor diaria" column is appended in H.
$ws.Range("D1:D2").EntireColumn.Delete()

# Add the new trailing column header + value.
$ws.Range("H1").Value = "Valor diaria"
$ws.Range("H2").Value = 60

# Update the purchase value and registration date for the existing record.
$ws.Range("E2").Value = 80000
$ws.Range("G2").Value = "2024-05-18 23:22:37"
